# Cambio a còdigo climatico
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 7 new rows below the existing data (rows 9-15) to host the
# additional weather stations, shifting nothing else around.
$ws.Rows.Item(9).Resize(7).Insert()

# --- Identificador (A) and Modbus_Ala (D): simple "+1 from the row above" fill ---
$ws.Range("A9:A15").Formula = "=A8+1"
$ws.Range("D9:D15").Formula = "=D8+1"

# --- Modbus_Add (C): step changed from +4 to +2 per row, from row 3 down ---
$ws.Range("C3:C15").Formula = "=C2+2"

# --- New station names (Ciudad) ---
$ws.Range("B9").Value = "Esmeraldas,EC"
$ws.Range("B10").Value = "Manta,EC"
$ws.Range("B11").Value = "Santo Domingo,EC"
$ws.Range("B12").Value = "Machala,EC"
$ws.Range("B13").Value = "Puyo,EC"
$ws.Range("B14").Value = "Tena,EC"
$ws.Range("B15").Value = "Macas,EC"

# --- Activa (E) flag for the new rows ---
$ws.Range("E9:E15").Value = "x"

# --- Column B widened (best fit) to fit the longer city names ---
$ws.Columns.Item(2).ColumnWidth = 16.42

# --- Selection moves down to the freshly filled Activa column ---
$ws.Range("E9:E15").Select() | Out-Null

$wb.Save() | Out-Null
